# reviewdb.xlsx update ("Add files via upload")
#
# Semantic changes applied (per the canonical-OOXML diff):
#   1. G5 review-status cell changes from "confirm" to "no".
#   2. Column B (appid/keyword "keyword" column, rows 2-11) is restyled to
#      match column A's font (Mangal 10pt) instead of the default Arial.
#   3. Column A width is narrowed and column F width is narrowed.
#   4. The sheet's active selection moves from G12 to D2:D11 (with the
#      window scrolled back to show column A, i.e. the default top-left).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the review-status value in G5: "confirm" -> "no"
$ws.Range("G5").Value = "no"

# 2. Re-style column B (rows 2-11) with the Mangal font used elsewhere in
#    the sheet (same font as column A / F / G), picking up the workbook's
#    existing "style 1" cell format.
$ws.Range("B2:B11").Font.Name = "Mangal"

# 3. Narrow columns A and F.
$ws.Columns.Item(1).ColumnWidth = 22.25
$ws.Columns.Item(6).ColumnWidth = 32.751

# 4. Move the selection to D2:D11 (active cell D2), scrolling the view back
#    to the default top-left (column A).
$ws.Activate()
$ws.Range("D2:D11").Select()
